$wb = $excel.ActiveWorkbook

# Values to write into C11:I11 (GDP, UEMP, CPI, LTRate, EURUSD, WTI, RPP)
$columns = @("C", "D", "E", "F", "G", "H", "I")
$values = @(
    0.6876068028317803,
    -0.3000000000000007,
    0.5459003767748243,
    -0.019000000000000017,
    1.5829618029997903,
    16.12947350163202,
    0.529961178858547
)

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])11"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
